$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data block: columns B-E for rows 29-32 (growth re-derivation from the
# "Previous measurements" block in rows 14-20), mirroring the existing G:J block. ---

# Row 29 (T0 / "1-2")
$ws.Range("B29").Value = 145
$ws.Range("C29").Formula = "=B16"
$ws.Range("D29:E29").Formula = "=C16"

# Row 30 (T1 / "2-1")
$ws.Range("B30").Value = 235
$ws.Range("C30").Formula = "=10*B17"
$ws.Range("D30:E30").Formula = "=10*C17"

# Row 31 (T2 / "3-2")
$ws.Range("B31").Value = 290
$ws.Range("C31").Formula = "=10*6*B18"
$ws.Range("D31:E31").Formula = "=10*6*C18"

# Row 32 (T3 / "4-2")
$ws.Range("B32").Value = 340
$ws.Range("C32").Formula = "=10*6*6*B20"
$ws.Range("D32:E32").Formula = "=10*6*6*C20"

# New 3-decimal number format applied to the first new row block (C29:E29)
$ws.Range("C29:E29").NumberFormat = "0.000"

# --- Column widths: drop the old column E width, give C/D/E their own widths ---
$ws.Columns("C").ColumnWidth = 9.625
$ws.Columns("D").ColumnWidth = 9.5
$ws.Columns("E").ColumnWidth = 10.375

# --- Page setup: explicit portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- View state: selection moves to C35 ---
$ws.Range("C35").Select()
